# Implemented doMove and timing test
#
# On the "Bitboard" sheet, duplicate the little 8x8 "square index" board
# (B2:I9) into a new block (T2:AA9), then sort that new block's rows in
# descending order on its first column (like Data > Sort Descending on
# column T, no header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bitboard")

# 1) Copy the existing board (values + formatting) into the new location.
$ws.Range("B2:I9").Copy($ws.Range("T2"))

# 2) Re-order the copied rows (descending by the first column, i.e. T)
#    the same way Data > Sort would, filling in the final values directly
#    so each destination cell keeps the formatting it already has.
for ($i = 0; $i -lt 8; $i++) {
    $srcRow = 9 - $i      # 9,8,7,6,5,4,3,2
    $dstRow = 2 + $i      # 2,3,4,5,6,7,8,9
    for ($col = 20; $col -le 27; $col++) {
        $srcCol = $col - 18   # same column offset as B:I (20-18=2 -> B)
        $value = $ws.Cells.Item($srcRow, $srcCol).Value()
        $ws.Cells.Item($dstRow, $col).Value = $value
    }
}

# 3) Apply an actual Sort on the new range so Excel records the sort
#    state (sortState/sortCondition) against T2:AA9, sorted descending
#    on T2:T9, with no header row.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("T2:T9"), $null, 2, $null, $null)
$sortObj.SetRange($ws.Range("T2:AA9"))
$sortObj.Header = 2
$sortObj.Apply()

# 4) Update the active selection shown when the sheet is reopened.
$ws.Activate()
$ws.Range("E15").Select()
